$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New email/password data set (replaces the old rows 2-4 and extends through row 25)
$data = @(
    @("acdf@hotmail.com","ncwchewcdc"),
    @("acd2f@hotmail.com","xcvb4567"),
    @("acdf3@hotmail.com","ncwchewcdc"),
    @("acdf4@hotmail.com","xcvb4567"),
    @("acdf5@hotmail.com","ncwchewcdc"),
    @("acdf6@hotmail.com","xcvb4567"),
    @("acdf7@hotmail.com","ncwchewcdc"),
    @("acdf8@hotmail.com","xcvb4567"),
    @("acdf9@hotmail.com","ncwchewcdc"),
    @("acdf10@hotmail.com","xcvb4567"),
    @("acdf11@hotmail.com","ncwchewcdc"),
    @("acdf12@hotmail.com","xcvb4567"),
    @("acdf13@hotmail.com","ncwchewcdc"),
    @("acdf14@hotmail.com","xcvb4567"),
    @("acdf15@hotmail.com","ncwchewcdc"),
    @("acdf16@hotmail.com","xcvb4567"),
    @("acdf17@hotmail.com","ncwchewcdc"),
    @("acdf18@hotmail.com","xcvb4567"),
    @("acdf19@hotmail.com","ncwchewcdc"),
    @("acdf20@hotmail.com","xcvb4567"),
    @("acdf21@hotmail.com","ncwchewcdc"),
    @("acdf22@hotmail.com","xcvb4567"),
    @("acdf23@hotmail.com","ncwchewcdc"),
    @("acdf24@hotmail.com","xcvb4567")
)

$row = 2
foreach ($pair in $data) {
    $email = $pair[0]
    $pw = $pair[1]
    $ws.Cells.Item($row, 1).Value = $email
    $ws.Cells.Item($row, 2).Value = $pw
    $ws.Hyperlinks.Add($ws.Cells.Item($row, 1), "mailto:$email")
    $row = $row + 1
}

# Center align every populated cell in columns A and B (header stays as-is)
$ws.Range("A2:B25").HorizontalAlignment = -4108

# Move the cursor roughly where the author left it
$ws.Range("B24").Select()
